# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Rate laws" sheet: change the rate-law expression for reaction_1
#    from "k1 * S1[c] * compartment" to "k1 * S1[c]" and move the
#    selection from C2 to C3.
# ---------------------------------------------------------------------
$rateLaws = $wb.Worksheets.Item("Rate laws")
$rateLaws.Range("C2").Value = "k1 * S1[c]"
$rateLaws.Activate()
$rateLaws.Range("C3").Select()

# ---------------------------------------------------------------------
# 2. "Parameters" sheet: add a "Submodels" value (ode_submodel) to the
#    new D column for both parameter rows, widen column G, and zoom to
#    130%.
# ---------------------------------------------------------------------
$params = $wb.Worksheets.Item("Parameters")
$params.Cells.Item(2, 4).Value = "ode_submodel"
$params.Cells.Item(3, 4).Value = "ode_submodel"
$params.Columns.Item(7).ColumnWidth = 35.8
$params.Activate()
$params.Range("A3").Select()
$excel.ActiveWindow.Zoom = 130

# ---------------------------------------------------------------------
# 3. "Submodels" sheet: move the selection from D2 to A2.
# ---------------------------------------------------------------------
$submodels = $wb.Worksheets.Item("Submodels")
$submodels.Activate()
$submodels.Range("A2").Select()

# ---------------------------------------------------------------------
# 4. "References" sheet becomes the active/selected sheet & tab.
# ---------------------------------------------------------------------
$references = $wb.Worksheets.Item("References")
$references.Activate()

Write-Host "edits applied"
